$d = $word.ActiveDocument

# The page footer ("Ver no Jupiter Salvar em pdf Salvar em docx" and the
# following "© 2020 ..." copyright notice), together with the blank
# paragraph that separates them from the preceding "Requisitos" text, is
# being removed from the end of the document. Locate the "Ver no Jupiter"
# paragraph by its text so the edit is resilient to the exact paragraph
# index, then delete it, the blank paragraph right before it, and the
# copyright paragraph right after it — leaving the blank paragraph (and
# the page-break paragraph) that follow untouched.

$count = $d.Paragraphs.Count
$jupiterIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Ver no Jupiter*") {
        $jupiterIndex = $i
        break
    }
}

if ($jupiterIndex -eq -1) {
    throw "Could not find the 'Ver no Jupiter...' paragraph"
}

$blankIndex = $jupiterIndex - 1
$copyrightIndex = $jupiterIndex + 1

$start = $d.Paragraphs.Item($blankIndex).Range.Start
$end = $d.Paragraphs.Item($copyrightIndex + 1).Range.Start
$d.Range($start, $end).Delete()
